# Apply updated loading-percent results for the 380 kV case (rows 2-25)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.13865601605993
$ws.Range("C2").Value = 11.65082694929254
$ws.Range("D2").Value = 6.023976823040961
$ws.Range("E2").Value = 16.23622155853771
$ws.Range("G2").Value = 3.697928491171754
$ws.Range("I2").Value = 29.90177494866774
$ws.Range("K2").Value = 15.16543590486079
$ws.Range("N2").Value = 22.06050047206049
$ws.Range("B3").Value = 13.85170263140027
$ws.Range("C3").Value = 11.32067020127965
$ws.Range("D3").Value = 5.913192090420042
$ws.Range("E3").Value = 15.32941049253146
$ws.Range("G3").Value = 3.701853157339202
$ws.Range("I3").Value = 29.76270398321813
$ws.Range("K3").Value = 14.93783338704886
$ws.Range("N3").Value = 22.05950917582567
$ws.Range("B4").Value = 13.67751017626312
$ws.Range("C4").Value = 11.11709939598795
$ws.Range("D4").Value = 5.846089042246787
$ws.Range("E4").Value = 14.7501427834372
$ws.Range("G4").Value = 3.704383117743995
$ws.Range("I4").Value = 29.68253631718972
$ws.Range("K4").Value = 14.80134307409016
$ws.Range("N4").Value = 22.06061651814846
$ws.Range("B5").Value = 13.60713763964477
$ws.Range("C5").Value = 11.0340674721524
$ws.Range("D5").Value = 5.819014409959433
$ws.Range("E5").Value = 14.50871092442595
$ws.Range("G5").Value = 3.705444452881163
$ws.Range("I5").Value = 29.65119505850527
$ws.Range("K5").Value = 14.74661331460048
$ws.Range("N5").Value = 22.06149824485174
$ws.Range("B6").Value = 13.59549253296952
$ws.Range("C6").Value = 11.02027973861643
$ws.Range("D6").Value = 5.814536163320486
$ws.Range("E6").Value = 14.46830583181273
$ws.Range("G6").Value = 3.705622524009911
$ws.Range("I6").Value = 29.6460714453653
$ws.Range("K6").Value = 14.73758136329682
$ws.Range("N6").Value = 22.0616706231403
$ws.Range("B7").Value = 13.67655848306256
$ws.Range("C7").Value = 11.1159797140321
$ws.Range("D7").Value = 5.845722757647496
$ws.Range("E7").Value = 14.74690810840996
$ws.Range("G7").Value = 3.704397308202104
$ws.Range("I7").Value = 29.68210824460098
$ws.Range("K7").Value = 14.80060126783503
$ws.Range("N7").Value = 22.06062666800201
$ws.Range("B8").Value = 14.03936548969023
$ws.Range("C8").Value = 11.53724759478585
$ws.Range("D8").Value = 5.985611037609724
$ws.Range("E8").Value = 15.92837027247081
$ws.Range("G8").Value = 3.699256847513378
$ws.Range("I8").Value = 29.85274856129362
$ws.Range("K8").Value = 15.08632581637091
$ws.Range("N8").Value = 22.05980191257427
$ws.Range("B9").Value = 14.76191217656282
$ws.Range("C9").Value = 12.35081179888962
$ws.Range("D9").Value = 6.265504203119374
$ws.Range("E9").Value = 18.05704931936777
$ws.Range("G9").Value = 3.690124177728708
$ws.Range("I9").Value = 30.22810159358987
$ws.Range("K9").Value = 15.66926245505778
$ws.Range("N9").Value = 22.07184644977569
$ws.Range("B10").Value = 15.29330057379214
$ws.Range("C10").Value = 12.93373019518559
$ws.Range("D10").Value = 6.472269970564856
$ws.Range("E10").Value = 19.67440724929947
$ws.Range("G10").Value = 3.683983804377266
$ws.Range("I10").Value = 30.52765013848485
$ws.Range("K10").Value = 16.10698329077541
$ws.Range("N10").Value = 22.08908787217015
$ws.Range("B11").Value = 15.5339118673344
$ws.Range("C11").Value = 13.19435999189055
$ws.Range("D11").Value = 6.566116332805888
$ws.Range("E11").Value = 20.37030312630491
$ws.Range("G11").Value = 3.681312229421285
$ws.Range("I11").Value = 30.66882935666803
$ws.Range("K11").Value = 16.30724418952456
$ws.Range("N11").Value = 22.09876623412211
$ws.Range("B12").Value = 15.62476650850926
$ws.Range("C12").Value = 13.29230124649464
$ws.Range("D12").Value = 6.601586345293079
$ws.Range("E12").Value = 20.62800502361628
$ws.Range("G12").Value = 3.680317935948868
$ws.Range("I12").Value = 30.72297274879811
$ws.Range("K12").Value = 16.38316663682599
$ws.Range("N12").Value = 22.10269587853317
$ws.Range("B13").Value = 15.60521239494643
$ws.Range("C13").Value = 13.27124286178113
$ws.Range("D13").Value = 6.59395081984748
$ws.Range("E13").Value = 20.57276252522426
$ws.Range("G13").Value = 3.680531304163464
$ws.Range("I13").Value = 30.71128207854387
$ws.Range("K13").Value = 16.36681265170269
$ws.Range("N13").Value = 22.1018377770847
$ws.Range("B14").Value = 15.54139228214685
$ws.Range("C14").Value = 13.20243332324437
$ws.Range("D14").Value = 6.569036032229804
$ws.Range("E14").Value = 20.39162082624294
$ws.Range("G14").Value = 3.681230080831815
$ws.Range("I14").Value = 30.67327024553844
$ws.Range("K14").Value = 16.31348903867758
$ws.Range("N14").Value = 22.09908422006371
$ws.Range("B15").Value = 15.50226397497379
$ws.Range("C15").Value = 13.16018453203501
$ws.Range("D15").Value = 6.553765136318446
$ws.Range("E15").Value = 20.2799096026816
$ws.Range("G15").Value = 3.681660360523758
$ws.Range("I15").Value = 30.6500749500633
$ws.Range("K15").Value = 16.28083602844191
$ws.Range("N15").Value = 22.09743207470302
$ws.Range("B16").Value = 15.27754497365452
$ws.Range("C16").Value = 12.91659722560956
$ws.Range("D16").Value = 6.466129412603956
$ws.Range("E16").Value = 19.6281125954087
$ws.Range("G16").Value = 3.684160836290265
$ws.Range("I16").Value = 30.51852053386552
$ws.Range("K16").Value = 16.09391216903469
$ws.Range("N16").Value = 22.08849236360238
$ws.Range("B17").Value = 15.13932846638173
$ws.Range("C17").Value = 12.76592648238887
$ws.Range("D17").Value = 6.412286461997397
$ws.Range("E17").Value = 19.2178535704621
$ws.Range("G17").Value = 3.685725879374928
$ws.Range("I17").Value = 30.4390574201486
$ws.Range("K17").Value = 15.9794776972881
$ws.Range("N17").Value = 22.08347873568638
$ws.Range("B18").Value = 15.05973014367199
$ws.Range("C18").Value = 12.67884248801823
$ws.Range("D18").Value = 6.381299541795227
$ws.Range("E18").Value = 18.97804883296994
$ws.Range("G18").Value = 3.686637514497831
$ws.Range("I18").Value = 30.39381633992925
$ws.Range("K18").Value = 15.9137697017089
$ws.Range("N18").Value = 22.08076772177998
$ws.Range("B19").Value = 15.03276546247245
$ws.Range("C19").Value = 12.64928804462821
$ws.Range("D19").Value = 6.370805987633801
$ws.Range("E19").Value = 18.89619536766716
$ws.Range("G19").Value = 3.686948151292039
$ws.Range("I19").Value = 30.37857890876465
$ws.Range("K19").Value = 15.8915436291916
$ws.Range("N19").Value = 22.0798794544782
$ws.Range("B20").Value = 15.15405291420986
$ws.Range("C20").Value = 12.78201008913657
$ws.Range("D20").Value = 6.418020259959074
$ws.Range("E20").Value = 19.26192294724207
$ws.Range("G20").Value = 3.68555809234177
$ws.Range("I20").Value = 30.44746855614466
$ws.Range("K20").Value = 15.99164843904328
$ws.Range("N20").Value = 22.0839945661657
$ws.Range("B21").Value = 15.56014560139189
$ws.Range("C21").Value = 13.22266555935702
$ws.Range("D21").Value = 6.576356229659108
$ws.Range("E21").Value = 20.44498416589598
$ws.Range("G21").Value = 3.681024362681265
$ws.Range("I21").Value = 30.68441693844069
$ws.Range("K21").Value = 16.32914968269593
$ws.Range("N21").Value = 22.09988581800175
$ws.Range("B22").Value = 15.82399548776945
$ws.Range("C22").Value = 13.50622673666987
$ws.Range("D22").Value = 6.679428167205748
$ws.Range("E22").Value = 21.18428353797158
$ws.Range("G22").Value = 3.67816251770814
$ws.Range("I22").Value = 30.84323868622187
$ws.Range("K22").Value = 16.55020550446523
$ws.Range("N22").Value = 22.11181460616646
$ws.Range("B23").Value = 15.6833472827542
$ws.Range("C23").Value = 13.35532118435977
$ws.Range("D23").Value = 6.62446601044408
$ws.Range("E23").Value = 20.79279542635584
$ws.Range("G23").Value = 3.679680719853255
$ws.Range("I23").Value = 30.75811842876621
$ws.Range("K23").Value = 16.43220450777133
$ws.Range("N23").Value = 22.10530657478527
$ws.Range("B24").Value = 15.14739640682981
$ws.Range("C24").Value = 12.7747401198978
$ws.Range("D24").Value = 6.415428106405282
$ws.Range("E24").Value = 19.24201146119782
$ws.Range("G24").Value = 3.685633911874266
$ws.Range("I24").Value = 30.4436644988599
$ws.Range("K24").Value = 15.98614578579718
$ws.Range("N24").Value = 22.08376082543087
$ws.Range("B25").Value = 14.56591270178041
$ws.Range("C25").Value = 12.13283644866716
$ws.Range("D25").Value = 6.189426586638576
$ws.Range("E25").Value = 17.50252676118611
$ws.Range("G25").Value = 3.692494209391998
$ws.Range("I25").Value = 30.12230068333216
$ws.Range("K25").Value = 15.50959004999789
$ws.Range("N25").Value = 22.06711862666802
